$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 15149.357
$ws.Range("I19").Value = 871.2857
$ws.Range("J19").Value = 29427.428
$ws.Range("K19").Value = 871.2857
$ws.Range("L19").Value = 29427.428
$ws.Range("M19").Value = -696.2857
$ws.Range("N19").Value = -29777.428
$ws.Range("H51").Value = 1940
$ws.Range("I51").Value = 1725
$ws.Range("J51").Value = 2083.3333
$ws.Range("K51").Value = 1725
$ws.Range("L51").Value = 2083.3333
$ws.Range("M51").Value = -1241
$ws.Range("N51").Value = -3051.3333
$ws.Range("H92").Value = 15874097
$ws.Range("I92").Value = 22223040
$ws.Range("J92").Value = 1739.6666
$ws.Range("K92").Value = 22223040
$ws.Range("L92").Value = 1739.6666
$ws.Range("M92").Value = -22221792
$ws.Range("N92").Value = -4235.6666
$ws.Range("H96").Value = 1150.5
$ws.Range("J96").Value = 2423.5
$ws.Range("L96").Value = 7270.5
$ws.Range("N96").Value = -10016.5
$ws.Range("H100").Value = 2640.16
$ws.Range("I100").Value = 2495.7778
$ws.Range("J100").Value = 3011.4285
$ws.Range("K100").Value = 2495.7778
$ws.Range("L100").Value = 3011.4285
$ws.Range("M100").Value = -1954.7778
$ws.Range("N100").Value = -4093.4285
$ws.Range("H115").Value = 4171.25
$ws.Range("I115").Value = 5342.5
$ws.Range("K115").Value = 16027.5
$ws.Range("M115").Value = -14460.5
$ws.Range("H132").Value = 4120.8286
$ws.Range("I132").Value = 3846.1292
$ws.Range("J132").Value = 6249.75
$ws.Range("K132").Value = 11538.3876
$ws.Range("L132").Value = 18749.25
$ws.Range("M132").Value = -9008.3876
$ws.Range("N132").Value = -23809.25
$ws.Range("H138").Value = 2437.158
$ws.Range("I138").Value = 4609.5
$ws.Range("J138").Value = 2082.4897
$ws.Range("K138").Value = 13828.5
$ws.Range("L138").Value = 6247.4691
$ws.Range("M138").Value = -8688.5
$ws.Range("N138").Value = -16527.4691
$ws.Range("H141").Value = 5070.3076
$ws.Range("I141").Value = 2197.0527
$ws.Range("K141").Value = 6591.158100000001
$ws.Range("M141").Value = -1411.158100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4137.4546
$ws.Range("I45").Value = 3502.4
$ws.Range("J45").Value = 4666.6665
$ws.Range("K45").Value = 3502.4
$ws.Range("L45").Value = 4666.6665
$ws.Range("M45").Value = -3125.4
$ws.Range("N45").Value = -5420.6665
$ws.Range("H96").Value = 34853.54
$ws.Range("J96").Value = 34853.54
$ws.Range("L96").Value = 34853.54
$ws.Range("N96").Value = -40345.54
$ws.Range("H122").Value = 58252.945
$ws.Range("I122").Value = 68943.53
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 206830.59
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -204380.59
$ws.Range("N122").Value = -19300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 50001970
$ws.Range("I86").Value = 62501990
$ws.Range("J86").Value = 1876.75
$ws.Range("K86").Value = 62501990
$ws.Range("L86").Value = 1876.75
$ws.Range("M86").Value = -62500867
$ws.Range("N86").Value = -4122.75
$ws.Range("H89").Value = 50001970
$ws.Range("I89").Value = 62501990
$ws.Range("J89").Value = 1876.75
$ws.Range("K89").Value = 312509950
$ws.Range("L89").Value = 9383.75
$ws.Range("M89").Value = -312504334
$ws.Range("N89").Value = -20615.75
$ws.Range("H99").Value = 1807
$ws.Range("I99").Value = 2021.6666
$ws.Range("J99").Value = 1485
$ws.Range("K99").Value = 2021.6666
$ws.Range("L99").Value = 1485
$ws.Range("M99").Value = -523.6666
$ws.Range("N99").Value = -4481

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1810.04
$ws.Range("I134").Value = 1408.0588
$ws.Range("J134").Value = 2664.25
$ws.Range("K134").Value = 4224.1764
$ws.Range("L134").Value = 7992.75
$ws.Range("M134").Value = -1689.1764
$ws.Range("N134").Value = -13062.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 983.5
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 3000
$ws.Range("N5").Value = -3224
$ws.Range("H23").Value = 62500130
$ws.Range("I23").Value = 111.6
$ws.Range("J23").Value = 90909224
$ws.Range("K23").Value = 334.8
$ws.Range("L23").Value = 272727672
$ws.Range("M23").Value = -99.79999999999995
$ws.Range("N23").Value = -272728142
$ws.Range("H110").Value = 12502.143
$ws.Range("J110").Value = 14148.117
$ws.Range("L110").Value = 42444.351
$ws.Range("N110").Value = -50624.351
$ws.Range("H113").Value = 1076.0952
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 1199.8823
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 3599.6469
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -7939.6469
$ws.Range("H131").Value = 846.2105
$ws.Range("I131").Value = 244
$ws.Range("J131").Value = 1061.2858
$ws.Range("K131").Value = 732
$ws.Range("L131").Value = 3183.8574
$ws.Range("M131").Value = 4308
$ws.Range("N131").Value = -13263.8574
$ws.Range("H134").Value = 4297.8667
$ws.Range("I134").Value = 2408.1177
$ws.Range("J134").Value = 6769.077
$ws.Range("K134").Value = 7224.353099999999
$ws.Range("L134").Value = 20307.231
$ws.Range("M134").Value = -2154.353099999999
$ws.Range("N134").Value = -30447.231
$ws.Range("H135").Value = 983.5
$ws.Range("J135").Value = 1000
$ws.Range("L135").Value = 9000
$ws.Range("N135").Value = -14070
$ws.Range("H139").Value = 2624.5789
$ws.Range("I139").Value = 2619.818
$ws.Range("J139").Value = 2631.125
$ws.Range("K139").Value = 7859.454000000001
$ws.Range("L139").Value = 7893.375
$ws.Range("M139").Value = -2719.454000000001
$ws.Range("N139").Value = -18173.375
$ws.Range("H140").Value = 1189.7727
$ws.Range("I140").Value = 1006.7778
$ws.Range("K140").Value = 3020.3334
$ws.Range("M140").Value = 2159.6666
$ws.Range("H141").Value = 5551.5
$ws.Range("I141").Value = 4470
$ws.Range("J141").Value = 7281.9
$ws.Range("K141").Value = 13410
$ws.Range("L141").Value = 21845.7
$ws.Range("M141").Value = -8230
$ws.Range("N141").Value = -32205.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5739.778
$ws.Range("I70").Value = 5617.0884
$ws.Range("K70").Value = 5617.0884
$ws.Range("M70").Value = -5347.0884
$ws.Range("H73").Value = 5739.778
$ws.Range("I73").Value = 5617.0884
$ws.Range("K73").Value = 5617.0884
$ws.Range("M73").Value = -4681.0884
$ws.Range("H97").Value = 2652.818
$ws.Range("I97").Value = 2025.5714
$ws.Range("J97").Value = 3750.5
$ws.Range("K97").Value = 2025.5714
$ws.Range("L97").Value = 3750.5
$ws.Range("M97").Value = -1529.5714
$ws.Range("N97").Value = -4742.5
$ws.Range("H122").Value = 3295.9375
$ws.Range("I122").Value = 2671.923
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 8015.768999999999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -5565.768999999999
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 3177.05
$ws.Range("I132").Value = 3089.5
$ws.Range("J132").Value = 3235.4167
$ws.Range("K132").Value = 9268.5
$ws.Range("L132").Value = 9706.250100000001
$ws.Range("M132").Value = -6738.5
$ws.Range("N132").Value = -14766.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2991.4285
$ws.Range("I61").Value = 2563.3333
$ws.Range("K61").Value = 2563.3333
$ws.Range("M61").Value = -2361.3333
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H93").Value = 7369.3
$ws.Range("I93").Value = 10626.182
$ws.Range("J93").Value = 3388.6667
$ws.Range("K93").Value = 10626.182
$ws.Range("L93").Value = 3388.6667
$ws.Range("M93").Value = -9378.182000000001
$ws.Range("N93").Value = -5884.6667
$ws.Range("H100").Value = 1921.9445
$ws.Range("I100").Value = 1470.6
$ws.Range("J100").Value = 2486.125
$ws.Range("K100").Value = 1470.6
$ws.Range("L100").Value = 2486.125
$ws.Range("M100").Value = -929.5999999999999
$ws.Range("N100").Value = -3568.125
$ws.Range("H113").Value = 2991.4285
$ws.Range("I113").Value = 2563.3333
$ws.Range("K113").Value = 2563.3333
$ws.Range("M113").Value = -393.3332999999998
$ws.Range("H122").Value = 3114.1428
$ws.Range("I122").Value = 2749.75
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 8249.25
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -5799.25
$ws.Range("N122").Value = -15700
$ws.Range("H132").Value = 3413
$ws.Range("I132").Value = 2680.647
$ws.Range("J132").Value = 4544.8184
$ws.Range("K132").Value = 8041.941
$ws.Range("L132").Value = 13634.4552
$ws.Range("M132").Value = -5511.941
$ws.Range("N132").Value = -18694.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 32000
$ws.Range("J99").Value = 32000
$ws.Range("L99").Value = 32000
$ws.Range("N99").Value = -37990
$ws.Range("H122").Value = 4547.8096
$ws.Range("I122").Value = 2667.3333
$ws.Range("J122").Value = 5300
$ws.Range("K122").Value = 8001.999899999999
$ws.Range("L122").Value = 15900
$ws.Range("M122").Value = -5551.999899999999
$ws.Range("N122").Value = -20800
$ws.Range("H123").Value = 24666
$ws.Range("J123").Value = 24666
$ws.Range("L123").Value = 24666
$ws.Range("N123").Value = -34466
$ws.Range("H132").Value = 3706037.2
$ws.Range("J132").Value = 5378221.5
$ws.Range("L132").Value = 16134664.5
$ws.Range("N132").Value = -16139724.5
